$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 225, shifting rows 225:266 down to 226:267
$ws.Rows.Item(225).Insert()

# Populate the new row 225 with data
$ws.Cells.Item(225, 1).Value = 11
$ws.Cells.Item(225, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(225, 3).Value = 'Bíobío'
$ws.Cells.Item(225, 4).Value = 45015
$ws.Cells.Item(225, 5).Value = 8
$ws.Cells.Item(225, 6).Value = 100112003
$ws.Cells.Item(225, 7).Value = 'Ajo'
$ws.Cells.Item(225, 8).Value = 'Chino'
$ws.Cells.Item(225, 9).Value = 'Primera'
$ws.Cells.Item(225, 10).Value = 250
$ws.Cells.Item(225, 11).Value = 15000
$ws.Cells.Item(225, 12).Value = 16000
$ws.Cells.Item(225, 13).Value = 15480
$ws.Cells.Item(225, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(225, 15).Value = 'China'
$ws.Cells.Item(225, 16).Value = 1548
$ws.Cells.Item(225, 17).Value = 10
$ws.Cells.Item(225, 18).Value = 'Hortaliza'
